$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.08 = 15672.77 pesos`n✅ 15672.77 pesos = 4.06 = 972.76 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 245.33
$ws2.Range("O10").Value = 3845
$ws2.Range("N12").Value = 3860
$ws2.Range("O12").Value = 239.577
